$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (laura) updates ---------------------------------------------
$ws.Range("C2").Value = "'2024-09-09"
$ws.Range("D2").Value = 14
$ws.Range("E2").Value = 1182.6429000000001
$ws.Range("F2").Value = 14
$ws.Range("G2").Value = 14

# --- Row 3 (rocio) updates ----------------------------------------------
$ws.Range("C3").Value = "'2024-09-09"
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 1443.6667
$ws.Range("F3").Value = 14
$ws.Range("G3").Value = 12

# --- Row 4: global average wait time ------------------------------------
$ws.Range("E4").Value = 1303.1153999999999

# --- Row 6-7: global peak-hour indicators --------------------------------
$ws.Range("A6").Value = "Hora Pico Global"
$ws.Range("B6").Value = 14
$ws.Range("A7").Value = "Turnos en Hora Pico Global"
$ws.Range("B7").Value = 26

# --- Row 9-10: totals summary --------------------------------------------
$ws.Range("A9").Value = "Total Clientes Atendidos"
$ws.Range("B9").Value = 26
$ws.Range("A10").Value = "Tiempo Promedio de Espera (min)"
$ws.Range("B10").Value = 1303.1153999999999

# --- Column widths (best-effort match to authored layout) ---------------
$ws.Columns.Item(1).ColumnWidth = 47.833333333333336
$ws.Columns.Item(2).ColumnWidth = 23
$ws.Columns.Item(3).ColumnWidth = 16.666666666666668
$ws.Columns.Item(5).ColumnWidth = 20.333333333333332
$ws.Columns.Item(7).ColumnWidth = 19.166666666666668

# --- Selection / view state ----------------------------------------------
$ws.Range("C23").Select()
